# Regenerate the "K" column (column G) values for each data row.
# These new values replace the previous (Strike#-derived) values with
# freshly calculated counts (s_vals), per the commit message:
# "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 2
    6  = 2
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 2
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 3
    22 = 0
    23 = 1
    24 = 2
    25 = 3
    26 = 0
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 3
    33 = 0
    34 = 1
    35 = 1
    36 = 0
    37 = 0
    38 = 3
    39 = 1
    40 = 2
    41 = 2
    42 = 1
    43 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
